$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.195784431877428
$ws.Range("C2").Value = 0.3505742835032208
$ws.Range("E2").Value = 0.4295085890014434
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.0974250560186789
$ws.Range("H2").Value = 0.2682436956180538
$ws.Range("I2").Value = 0.1596043950206705
$ws.Range("O2").Value = 0.6201320040570124

$ws.Range("B3").Value = 1.043962204355239
$ws.Range("C3").Value = 0.3116550348671012
$ws.Range("E3").Value = 0.3746495323167949
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.09910359483972542
$ws.Range("H3").Value = 0.2734605452509911
$ws.Range("I3").Value = 0.1660391165689266
$ws.Range("O3").Value = 0.634349020292035

$ws.Range("B4").Value = 0.9503500411424284
$ws.Range("C4").Value = 0.287653993403552
$ws.Range("E4").Value = 0.3410483920271048
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.1003810139088657
$ws.Range("H4").Value = 0.276917597517091
$ws.Range("I4").Value = 0.1702482771225904
$ws.Range("O4").Value = 0.6441244567818103

$ws.Range("B5").Value = 0.9121063457178593
$ws.Range("C5").Value = 0.2778478352446712
$ws.Range("E5").Value = 0.3273747622888834
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.1009630262687011
$ws.Range("H5").Value = 0.2783900197801721
$ws.Range("I5").Value = 0.1720281635846057
$ws.Range("O5").Value = 0.6483693889936148

$ws.Range("B6").Value = 0.9057502944375528
$ws.Range("C6").Value = 0.2762180092985318
$ws.Range("E6").Value = 0.3251053725941802
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.1010633629369373
$ws.Range("H6").Value = 0.2786383526540455
$ws.Range("I6").Value = 0.1723276061204075
$ws.Range("O6").Value = 0.6490899926976894

$ws.Range("B7").Value = 0.9498346582125805
$ws.Range("C7").Value = 0.2875218465555065
$ws.Range("E7").Value = 0.3408639094686379
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.1003886151081872
$ws.Range("H7").Value = 0.2769371976825425
$ws.Range("I7").Value = 0.1702720200653975
$ws.Range("O7").Value = 0.6441806494275539

$ws.Range("B8").Value = 1.143519239471175
$ws.Range("C8").Value = 0.3371770321460019
$ws.Range("E8").Value = 0.4105750445756939
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.09795223237942707
$ws.Range("H8").Value = 0.2699896539270981
$ws.Range("I8").Value = 0.1617693081711322
$ws.Range("O8").Value = 0.624815888511435

$ws.Range("B9").Value = 1.520115483469397
$ws.Range("C9").Value = 0.4336945908939356
$ws.Range("E9").Value = 0.5480141449565252
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.09515976530711612
$ws.Range("H9").Value = 0.2583886878790551
$ws.Range("I9").Value = 0.1471579894610833
$ws.Range("O9").Value = 0.5952181273748778

$ws.Range("B10").Value = 1.794728392379511
$ws.Range("C10").Value = 0.5040546680731381
$ws.Range("E10").Value = 0.6495604917734283
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.09435576854083649
$ws.Range("H10").Value = 0.2511105906744504
$ws.Range("I10").Value = 0.1376987339919873
$ws.Range("O10").Value = 0.5786833671548237

$ws.Range("B11").Value = 1.919185743947537
$ws.Range("C11").Value = 0.5359378742575132
$ws.Range("E11").Value = 0.6959071633556988
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.0942689315619134
$ws.Range("H11").Value = 0.248072491237707
$ws.Range("I11").Value = 0.1336764023469681
$ws.Range("O11").Value = 0.5723150849017316

$ws.Range("B12").Value = 1.966245252617682
$ws.Range("C12").Value = 0.54799277516139
$ws.Range("E12").Value = 0.7134814030711993
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.0942768055061336
$ws.Range("H12").Value = 0.2469614763392372
$ws.Range("I12").Value = 0.1321939538289294
$ws.Range("O12").Value = 0.570071282091476

$ws.Range("B13").Value = 1.956113294489001
$ws.Range("C13").Value = 0.5453973747042937
$ws.Range("E13").Value = 0.7096953904892587
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.09427328665947954
$ws.Range("H13").Value = 0.2471989952185609
$ws.Range("I13").Value = 0.1325114089369535
$ws.Range("O13").Value = 0.5705470358152951

$ws.Range("B14").Value = 1.923058772988327
$ws.Range("C14").Value = 0.5369300148131515
$ws.Range("E14").Value = 0.6973525200568105
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.09426875879282903
$ws.Range("H14").Value = 0.2479802954854122
$ws.Range("I14").Value = 0.1335536222810561
$ws.Range("O14").Value = 0.5721271137265802

$ws.Range("B15").Value = 1.90280275122592
$ws.Range("C15").Value = 0.5317410708495913
$ws.Range("E15").Value = 0.689795303617899
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.09427131245272591
$ws.Range("H15").Value = 0.2484640084272769
$ws.Range("I15").Value = 0.1341973216182042
$ws.Range("O15").Value = 0.5731168553683403

$ws.Range("B16").Value = 1.786585200693821
$ws.Range("C16").Value = 0.5019684667934712
$ws.Range("E16").Value = 0.6465348544352167
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.09436711148332932
$ws.Range("H16").Value = 0.251314645330794
$ws.Range("I16").Value = 0.1379672819586979
$ws.Range("O16").Value = 0.5791229190759566

$ws.Range("B17").Value = 1.715168242037976
$ws.Range("C17").Value = 0.4836716536509016
$ws.Range("E17").Value = 0.6200363617091682
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.09449777869713927
$ws.Range("H17").Value = 0.2531334364681044
$ws.Range("I17").Value = 0.1403521829916037
$ws.Range("O17").Value = 0.5831041982574305

$ws.Range("B18").Value = 1.674047455840196
$ws.Range("C18").Value = 0.4731361888065635
$ws.Range("E18").Value = 0.6048093965388119
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.09459914800327596
$ws.Range("H18").Value = 0.2542052153210577
$ws.Range("I18").Value = 0.1417503176127333
$ws.Range("O18").Value = 0.5855025828962255

$ws.Range("B19").Value = 1.660117256599676
$ws.Range("C19").Value = 0.4695670884911465
$ws.Range("E19").Value = 0.5996562036081485
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.09463795217316573
$ws.Range("H19").Value = 0.2545725011459012
$ws.Range("I19").Value = 0.1422282259248242
$ws.Range("O19").Value = 0.5863332058234931

$ws.Range("B20").Value = 1.722775236657981
$ws.Range("C20").Value = 0.4856205889574312
$ws.Range("E20").Value = 0.622855681963415
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.09448115128672896
$ws.Range("H20").Value = 0.2529371658625479
$ws.Range("I20").Value = 0.1400955710812402
$ws.Range("O20").Value = 0.5826691457957196

$ws.Range("B21").Value = 1.932769603888687
$ws.Range("C21").Value = 0.5394175938148464
$ws.Range("E21").Value = 0.7009772610038283
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.09426897741103346
$ws.Range("H21").Value = 0.2477497362535885
$ws.Range("I21").Value = 0.1332463909850174
$ws.Range("O21").Value = 0.5716584396359679

$ws.Range("B22").Value = 2.069604974032814
$ws.Range("C22").Value = 0.5744684199427752
$ws.Range("E22").Value = 0.7521735745979896
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.09436817572044021
$ws.Range("H22").Value = 0.2445894866194749
$ws.Range("I22").Value = 0.129007574113815
$ws.Range("O22").Value = 0.565440814401299

$ws.Range("B23").Value = 1.996611637849242
$ws.Range("C23").Value = 0.5557713205332107
$ws.Range("E23").Value = 0.7248357891637198
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.0942932519197015
$ws.Range("H23").Value = 0.2462550478863648
$ws.Range("I23").Value = 0.1312480608928386
$ws.Range("O23").Value = 0.5686691295973674

$ws.Range("B24").Value = 1.719336308000948
$ws.Range("C24").Value = 0.4847395262848408
$ws.Range("E24").Value = 0.6215810444864616
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.09448858683060735
$ws.Range("H24").Value = 0.2530258184193315
$ws.Range("I24").Value = 0.140211501175191
$ws.Range("O24").Value = 0.5828654920809271

$ws.Range("B25").Value = 1.4185923225848
$ws.Range("C25").Value = 0.4076787669468445
$ws.Range("E25").Value = 0.5107419805693638
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.09569880217166116
$ws.Range("H25").Value = 0.2613092496933191
$ws.Range("I25").Value = 0.1508880435335684
$ws.Range("O25").Value = 0.6023175765766666
